$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of |S*|/n column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary rows 14-17: labels in column A, aggregate formulas in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"

$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style B14 (bold, size 12, vertically centered) then propagate the same
# format to B15:B17 via a format-only copy/paste so a single style entry
# is reused for all four cells.
$f = $ws.Range("B14").Font
$f.Bold = $true
$f.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# Make the new summary block the active selection, as in the authored file.
$ws.Range("A14:B17").Select()

# Match the page setup recorded by the authoring Excel session.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
